{"js": "// Update the date line and the 25 multiplication-fact answers in the table.\nconst replacements = [\n  [\"2025-02-09 Sunday\", \"2025-02-10 Monday\"],\n  [\"667\u00d76=4002\", \"455\u00d76=2730\"],\n  [\"544\u00d76=3264\", \"397\u00d77=2779\"],\n  [\"302\u00d78=2416\", \"778\u00d73=2334\"],\n  [\"550\u00d79=4950\", \"924\u00d72=1848\"],\n  [\"507\u00d77=3549\", \"683\u00d76=4098\"],\n  [\"954\u00d78=7632\", \"936\u00d79=8424\"],\n  [\"673\u00d77=4711\", \"278\u00d75=1390\"],\n  [\"963\u00d75=4815\", \"337\u00d77=2359\"],\n  [\"308\u00d75=1540\", \"298\u00d77=2086\"],\n  [\"859\u00d73=2577\", \"727\u00d79=6543\"],\n  [\"987\u00d73=2961\", \"927\u00d79=8343\"],\n  [\"575\u00d76=3450\", \"297\u00d78=2376\"],\n  [\"315\u00d73=945\", \"105\u00d76=630\"],\n  [\"279\u00d72=558\", \"741\u00d72=1482\"],\n  [\"351\u00d78=2808\", \"492\u00d76=2952\"],\n  [\"495\u00d79=4455\", \"834\u00d72=1668\"],\n  [\"450\u00d72=900\", \"589\u00d79=5301\"],\n  [\"237\u00d79=2133\", \"619\u00d75=3095\"],\n  [\"402\u00d79=3618\", \"977\u00d72=1954\"],\n  [\"137\u00d77=959\", \"607\u00d78=4856\"],\n  [\"163\u00d79=1467\", \"228\u00d79=2052\"],\n  [\"405\u00d73=1215\", \"416\u00d78=3328\"],\n  [\"301\u00d75=1505\", \"543\u00d73=1629\"],\n  [\"486\u00d79=4374\", \"682\u00d72=1364\"],\n  [\"825\u00d77=5775\", \"315\u00d75=1575\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const found = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  found.load(\"items\");\n  await context.sync();\n\n  for (const range of found.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Update the date line and the 25 multiplication-fact answers in the table.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2025-02-09 Sunday\", \"2025-02-10 Monday\"),\n    @(\"667\u00d76=4002\", \"455\u00d76=2730\"),\n    @(\"544\u00d76=3264\", \"397\u00d77=2779\"),\n    @(\"302\u00d78=2416\", \"778\u00d73=2334\"),\n    @(\"550\u00d79=4950\", \"924\u00d72=1848\"),\n    @(\"507\u00d77=3549\", \"683\u00d76=4098\"),\n    @(\"954\u00d78=7632\", \"936\u00d79=8424\"),\n    @(\"673\u00d77=4711\", \"278\u00d75=1390\"),\n    @(\"963\u00d75=4815\", \"337\u00d77=2359\"),\n    @(\"308\u00d75=1540\", \"298\u00d77=2086\"),\n    @(\"859\u00d73=2577\", \"727\u00d79=6543\"),\n    @(\"987\u00d73=2961\", \"927\u00d79=8343\"),\n    @(\"575\u00d76=3450\", \"297\u00d78=2376\"),\n    @(\"315\u00d73=945\", \"105\u00d76=630\"),\n    @(\"279\u00d72=558\", \"741\u00d72=1482\"),\n    @(\"351\u00d78=2808\", \"492\u00d76=2952\"),\n    @(\"495\u00d79=4455\", \"834\u00d72=1668\"),\n    @(\"450\u00d72=900\", \"589\u00d79=5301\"),\n    @(\"237\u00d79=2133\", \"619\u00d75=3095\"),\n    @(\"402\u00d79=3618\", \"977\u00d72=1954\"),\n    @(\"137\u00d77=959\", \"607\u00d78=4856\"),\n    @(\"163\u00d79=1467\", \"228\u00d79=2052\"),\n    @(\"405\u00d73=1215\", \"416\u00d78=3328\"),\n    @(\"301\u00d75=1505\", \"543\u00d73=1629\"),\n    @(\"486\u00d79=4374\", \"682\u00d72=1364\"),\n    @(\"825\u00d77=5775\", \"315\u00d75=1575\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Execute($null, $false, $false, $false, $false, $false, $true, 1, $false, $null, 2)\n}\n"}
